$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "49.453.43"
$ws.Range("D3").Value = "2.624.86"
$ws.Range("E3").Value = "  -1.47%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "111.09"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -2.68%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "325.05"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.15%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.523"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.73%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.999"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.543"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -2.79%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "39.28"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -5.14%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "20.05"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.85%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.0807"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -2.44%  "
$ws.Range("E13").Value = "  +1.30%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "7.34"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.74%  "
$ws.Range("D15").Value = "3.036.03"
$ws.Range("E15").Value = "  +0.49%  "
$ws.Range("D16").Value = "2.634.06"
$ws.Range("E16").Value = "  -1.55%  "
$ws.Range("E17").Value = "  -3.46%  "
$ws.Range("D18").Value = "49.365.16"
$ws.Range("E18").Value = "  -1.30%  "
$ws.Range("E19").Value = "  -2.36%  "
$ws.Range("E20").Value = "  -1.47%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.66"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.30%  "
$ws.Range("D22").Value = "0.0₃0945"
$ws.Range("E22").Value = "  -2.01%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "266.78"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -4.28%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "68.72"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -5.31%  "
$ws.Range("E25").Value = "  -2.61%  "
$ws.Range("E26").Value = "  +0.08%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "25.95"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -3.86%  "
$ws.Range("E28").Value = "  +1.27%  "
$ws.Range("E29").Value = "  -1.29%  "
$ws.Range("E30").Value = "  -2.25%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "34.48"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -5.27%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "49.56"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -1.90%  "
$ws.Range("E33").Value = "  +0.39%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0806"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +0.44%  "
$ws.Range("E35").Value = "  -0.20%  "
$ws.Range("E37").Value = "  +2.74%  "
$ws.Range("E38").Value = "  -3.56%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.08"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "128.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +2.26%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "22.63"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.59%  "
$ws.Range("B42").Value = "WEMIXToken"
$ws.Range("C42").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.24"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.24%  "
$ws.Range("B43").Value = "Stellar"
$ws.Range("C43").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.111"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.03%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0325"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +2.93%  "
$ws.Range("D45").Value = "2.036.23"
$ws.Range("E45").Value = "  -2.00%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.15"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +7.82%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.20"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -5.20%  "
$ws.Range("E48").Value = "  -3.96%  "
$ws.Range("E49").Value = "  -3.82%  "
$ws.Range("E50").Value = "  -4.21%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "58.32"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +0.67%  "
